# Update per-sheet "currentAveragePrice*" / computed profit columns (H-N)
# with freshly fetched market data, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2632580.2
$ws.Range("J19").Value = 1430.5
$ws.Range("L19").Value = 1430.5
$ws.Range("N19").Value = -1780.5
$ws.Range("H33").Value = 349.77777
$ws.Range("I33").Value = 269.8
$ws.Range("J33").Value = 449.75
$ws.Range("K33").Value = 269.8
$ws.Range("L33").Value = 449.75
$ws.Range("M33").Value = -40.80000000000001
$ws.Range("N33").Value = -907.75
$ws.Range("H55").Value = 184.66667
$ws.Range("I55").Value = 149.8
$ws.Range("J55").Value = 198.07692
$ws.Range("K55").Value = 149.8
$ws.Range("L55").Value = 198.07692
$ws.Range("M55").Value = 64.19999999999999
$ws.Range("N55").Value = -626.07692
$ws.Range("H125").Value = 1005.93335
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 1032.4166
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 9291.749400000001
$ws.Range("M125").Value = -5640
$ws.Range("N125").Value = -14211.7494
$ws.Range("H137").Value = 2464.0454
$ws.Range("I137").Value = 1916.3572
$ws.Range("J137").Value = 3422.5
$ws.Range("K137").Value = 5749.071599999999
$ws.Range("L137").Value = 10267.5
$ws.Range("M137").Value = -3199.071599999999
$ws.Range("N137").Value = -15367.5
$ws.Range("H141").Value = 140171.5
$ws.Range("I141").Value = 185632
$ws.Range("K141").Value = 556896
$ws.Range("M141").Value = -551716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4293.75
$ws.Range("I32").Value = 4124.2173
$ws.Range("J32").Value = 5073.6
$ws.Range("K32").Value = 4124.2173
$ws.Range("L32").Value = 5073.6
$ws.Range("M32").Value = -3837.2173
$ws.Range("N32").Value = -5647.6
$ws.Range("H132").Value = 2556.6924
$ws.Range("I132").Value = 1324.5333
$ws.Range("J132").Value = 4236.909
$ws.Range("K132").Value = 3973.5999
$ws.Range("L132").Value = 12710.727
$ws.Range("M132").Value = -1443.5999
$ws.Range("N132").Value = -17770.727

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 15562.667
$ws.Range("J81").Value = 15562.667
$ws.Range("L81").Value = 15562.667
$ws.Range("N81").Value = -17684.667
$ws.Range("H84").Value = 15562.667
$ws.Range("J84").Value = 15562.667
$ws.Range("L84").Value = 46688.001
$ws.Range("N84").Value = -57296.001
$ws.Range("H105").Value = 1600.4762
$ws.Range("I105").Value = 1596
$ws.Range("J105").Value = 1627.3334
$ws.Range("K105").Value = 1596
$ws.Range("L105").Value = 1627.3334
$ws.Range("M105").Value = 151
$ws.Range("N105").Value = -5121.3334
$ws.Range("H134").Value = 2212.5715
$ws.Range("I134").Value = 1577.64
$ws.Range("J134").Value = 3799.9
$ws.Range("K134").Value = 4732.92
$ws.Range("L134").Value = 11399.7
$ws.Range("M134").Value = -2197.92
$ws.Range("N134").Value = -16469.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22223402
$ws.Range("I16").Value = 37037936
$ws.Range("J16").Value = 1599.5
$ws.Range("K16").Value = 37037936
$ws.Range("L16").Value = 1599.5
$ws.Range("M16").Value = -37037649
$ws.Range("N16").Value = -2173.5
$ws.Range("H22").Value = 826.7143
$ws.Range("I22").Value = 260.2
$ws.Range("J22").Value = 1141.4445
$ws.Range("K22").Value = 260.2
$ws.Range("L22").Value = 1141.4445
$ws.Range("M22").Value = 89.80000000000001
$ws.Range("N22").Value = -1841.4445
$ws.Range("H31").Value = 2639.8
$ws.Range("I31").Value = 1140.9231
$ws.Range("J31").Value = 4263.5835
$ws.Range("K31").Value = 1140.9231
$ws.Range("L31").Value = 4263.5835
$ws.Range("M31").Value = -845.9231
$ws.Range("N31").Value = -4853.5835
$ws.Range("H34").Value = 2639.8
$ws.Range("I34").Value = 1140.9231
$ws.Range("J34").Value = 4263.5835
$ws.Range("K34").Value = 1140.9231
$ws.Range("L34").Value = 4263.5835
$ws.Range("M34").Value = -938.9231
$ws.Range("N34").Value = -4667.5835
$ws.Range("H113").Value = 22223402
$ws.Range("I113").Value = 37037936
$ws.Range("J113").Value = 1599.5
$ws.Range("K113").Value = 37037936
$ws.Range("L113").Value = 1599.5
$ws.Range("M113").Value = -37035766
$ws.Range("N113").Value = -5939.5
$ws.Range("H132").Value = 2601.0386
$ws.Range("I132").Value = 1858.9412
$ws.Range("J132").Value = 4002.7778
$ws.Range("K132").Value = 5576.8236
$ws.Range("L132").Value = 12008.3334
$ws.Range("M132").Value = -3046.8236
$ws.Range("N132").Value = -17068.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1945.6364
$ws.Range("J22").Value = 1945.6364
$ws.Range("L22").Value = 5836.9092
$ws.Range("N22").Value = -6174.9092
$ws.Range("H27").Value = 1945.6364
$ws.Range("J27").Value = 1945.6364
$ws.Range("L27").Value = 5836.9092
$ws.Range("N27").Value = -6040.9092
$ws.Range("H131").Value = 5682596.5
$ws.Range("J131").Value = 792.52325
$ws.Range("L131").Value = 2377.56975
$ws.Range("N131").Value = -12457.56975
$ws.Range("H132").Value = 2830.3635
$ws.Range("I132").Value = 568.75
$ws.Range("K132").Value = 5118.75
$ws.Range("M132").Value = -2588.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 27860
$ws.Range("J42").Value = 27860
$ws.Range("L42").Value = 27860
$ws.Range("N42").Value = -28830
$ws.Range("H63").Value = 16900
$ws.Range("J63").Value = 16900
$ws.Range("L63").Value = 16900
$ws.Range("N63").Value = -18272
$ws.Range("H66").Value = 16900
$ws.Range("J66").Value = 16900
$ws.Range("L66").Value = 50700
$ws.Range("N66").Value = -57564
$ws.Range("H70").Value = 6261.2593
$ws.Range("I70").Value = 5787
$ws.Range("J70").Value = 8348
$ws.Range("K70").Value = 5787
$ws.Range("L70").Value = 8348
$ws.Range("M70").Value = -5517
$ws.Range("N70").Value = -8888
$ws.Range("H73").Value = 6261.2593
$ws.Range("I73").Value = 5787
$ws.Range("J73").Value = 8348
$ws.Range("K73").Value = 5787
$ws.Range("L73").Value = 8348
$ws.Range("M73").Value = -4851
$ws.Range("N73").Value = -10220
$ws.Range("H115").Value = 27860
$ws.Range("J115").Value = 27860
$ws.Range("L115").Value = 27860
$ws.Range("N115").Value = -30210
$ws.Range("H132").Value = 4204.4
$ws.Range("I132").Value = 2435.1428
$ws.Range("K132").Value = 7305.428400000001
$ws.Range("M132").Value = -4775.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 684.3200000000001
$ws.Range("I68").Value = 686.1836499999999
$ws.Range("J68").Value = 593
$ws.Range("K68").Value = 686.1836499999999
$ws.Range("L68").Value = 593
$ws.Range("M68").Value = 62.81635000000006
$ws.Range("N68").Value = -2091
$ws.Range("H71").Value = 684.3200000000001
$ws.Range("I71").Value = 686.1836499999999
$ws.Range("J71").Value = 593
$ws.Range("K71").Value = 3430.91825
$ws.Range("L71").Value = 2965
$ws.Range("M71").Value = 313.0817500000003
$ws.Range("N71").Value = -10453
$ws.Range("H100").Value = 1450.625
$ws.Range("I100").Value = 1401.5
$ws.Range("J100").Value = 1467
$ws.Range("K100").Value = 1401.5
$ws.Range("L100").Value = 1467
$ws.Range("M100").Value = -860.5
$ws.Range("N100").Value = -2549

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 28900
$ws.Range("J64").Value = 28900
$ws.Range("L64").Value = 28900
$ws.Range("N64").Value = -29396
$ws.Range("H67").Value = 28900
$ws.Range("J67").Value = 28900
$ws.Range("L67").Value = 28900
$ws.Range("N67").Value = -30616
$ws.Range("H126").Value = 3962.5
$ws.Range("I126").Value = 2140
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 6420
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -3950
$ws.Range("N126").Value = -25940
$ws.Range("H132").Value = 15875352
$ws.Range("I132").Value = 1424.2
$ws.Range("K132").Value = 4272.6
$ws.Range("M132").Value = -1742.6
$ws.Range("H136").Value = 6162.316
$ws.Range("I136").Value = 4716.636
$ws.Range("J136").Value = 8150.125
$ws.Range("K136").Value = 14149.908
$ws.Range("L136").Value = 24450.375
$ws.Range("M136").Value = -11599.908
$ws.Range("N136").Value = -29550.375
